$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 59, pushing existing rows 59:135 down to 60:136.
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with a new weekly price record.
$ws.Cells.Item(59, 1).Value = 7
$ws.Cells.Item(59, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value = "Ñuble"
$ws.Cells.Item(59, 4).Value = 44413
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 5).Value = 16
$ws.Cells.Item(59, 6).Value = 100112009
$ws.Cells.Item(59, 7).Value = "Acelga"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 120
$ws.Cells.Item(59, 11).Value = 400
$ws.Cells.Item(59, 12).Value = 450
$ws.Cells.Item(59, 13).Value = 425
$ws.Cells.Item(59, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(59, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(59, 16).Value = 425
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"
